$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("F6").Value = "Ekran Alıntısı.JPG"
[void]$ws.Range("F7").Select()
